# Apply January 2022 regional table updates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated percentage strings (shared-string text edits) ---
# Row 2 = "Switzerland overall": Percentage Delta (95% CI)
$ws.Range("K2").Value = "6.2 (5.8-6.7)"

# Row 4 = "Region 2": Percentage Omicron (BA.1) (95% CI)
$ws.Range("S4").Value = "89.7 (88.2-91.1)"

# Row 5 = "Region 3": Percentage Omicron (BA.1) (95% CI)
$ws.Range("S5").Value = "91.7 (90.6-92.6)"

# --- Updated "Sequenced samples" (column B) and "Omicron (BA.1)" (column R) counts ---
# Row 2 = "Switzerland overall"
$ws.Range("B2").Value = 11287.0
$ws.Range("R2").Value = 10375.0

# Row 3 = "Region 1"
$ws.Range("B3").Value = 3316.0
$ws.Range("R3").Value = 3106.0

# Row 4 = "Region 2"
$ws.Range("B4").Value = 1655.0
$ws.Range("R4").Value = 1485.0

# Row 5 = "Region 3"
$ws.Range("B5").Value = 2925.0
$ws.Range("R5").Value = 2681.0

# Row 7 = "Region 5"
$ws.Range("B7").Value = 2042.0
$ws.Range("R7").Value = 1864.0
